# Update "want to go" (想去人数) counts on several rows across sheets.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 826
$ws1.Range("F4").Value  = 274
$ws1.Range("F5").Value  = 445
$ws1.Range("F8").Value  = 957
$ws1.Range("F9").Value  = 1025
$ws1.Range("F13").Value = 142
$ws1.Range("F16").Value = 25442
$ws1.Range("F17").Value = 2594
$ws1.Range("F22").Value = 135
$ws1.Range("F23").Value = 414
$ws1.Range("F24").Value = 228
$ws1.Range("F25").Value = 155
$ws1.Range("F28").Value = 104
$ws1.Range("F33").Value = 489

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F10").Value = 3833
$ws2.Range("F21").Value = 4166

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value = 196
$ws3.Range("F4").Value = 903

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 196
$ws4.Range("F4").Value  = 903
$ws4.Range("F5").Value  = 274
$ws4.Range("F6").Value  = 445
$ws4.Range("F15").Value = 957
$ws4.Range("F16").Value = 1025
$ws4.Range("F19").Value = 142
$ws4.Range("F28").Value = 2594
$ws4.Range("F35").Value = 414
$ws4.Range("F36").Value = 228
$ws4.Range("F40").Value = 104
$ws4.Range("F46").Value = 489
